# Auto-generated edit script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.569.72'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.637.08'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.91'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.503'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.93'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0844'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.863.54'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.16'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.71%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.606.96'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.62%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.528'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.39'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.54%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.575.44'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0₃0744'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '216.75'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.18%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.33'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.29'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('E23').Value = '  +16.97%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.39'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.28'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.91'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('E29').Value = '  +1.85%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0516'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.17'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.37'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.68%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.271.10'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +8.83%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.52'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.71%  '
$ws.Range('E36').Value = '  +1.10%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0176'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.39%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.513'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.802'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.28'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.93%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.801'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.37'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.773.67'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '93.56'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '55.28'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.62'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.18%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0961'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.37%  '
